$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 (A1) keeps displaying the same text "HK_G_acc_G"; the commit only
# changes which shared-string slot backs it (8 duplicate entries were
# inserted ahead of it upstream). Re-asserting the literal text reproduces
# the visible effect without needing to touch the shared-string pool by hand.
$ws.Range("A1").Value = "HK_G_acc_G"

# Rows 2-49: replace the computed Global-Threshold accuracy values with the
# freshly recalculated mean-based series from the commit.
$ws.Range("A2").Value = 48.141891891891895
$ws.Range("A3").Value = 48.310810810810814
$ws.Range("A4").Value = 48.47972972972973
$ws.Range("A5").Value = 47.80405405405405
$ws.Range("A6").Value = 48.310810810810814
$ws.Range("A7").Value = 48.310810810810814
$ws.Range("A8").Value = 48.986486486486484
$ws.Range("A9").Value = 49.1554054054054
$ws.Range("A10").Value = 48.64864864864865
$ws.Range("A11").Value = 48.47972972972973
$ws.Range("A12").Value = 48.986486486486484
$ws.Range("A13").Value = 49.83108108108108
$ws.Range("A14").Value = 49.66216216216216
$ws.Range("A15").Value = 50.16891891891891
$ws.Range("A16").Value = 50.16891891891891
$ws.Range("A17").Value = 50.50675675675676
$ws.Range("A18").Value = 51.01351351351351
$ws.Range("A19").Value = 50.8445945945946
$ws.Range("A20").Value = 48.986486486486484
$ws.Range("A21").Value = 49.1554054054054
$ws.Range("A22").Value = 48.986486486486484
$ws.Range("A23").Value = 47.2972972972973
$ws.Range("A24").Value = 46.95945945945946
$ws.Range("A25").Value = 46.95945945945946
$ws.Range("A26").Value = 48.47972972972973
$ws.Range("A27").Value = 48.47972972972973
$ws.Range("A28").Value = 49.32432432432432
$ws.Range("A29").Value = 50.16891891891891
$ws.Range("A30").Value = 50
$ws.Range("A31").Value = 49.49324324324324
$ws.Range("A32").Value = 47.46621621621622
$ws.Range("A33").Value = 46.95945945945946
$ws.Range("A34").Value = 47.80405405405405
$ws.Range("A35").Value = 47.63513513513514
$ws.Range("A36").Value = 48.310810810810814
$ws.Range("A37").Value = 51.520270270270274
$ws.Range("A38").Value = 47.46621621621622
$ws.Range("A39").Value = 47.63513513513514
$ws.Range("A40").Value = 47.97297297297297
$ws.Range("A41").Value = 49.1554054054054
$ws.Range("A42").Value = 48.47972972972973
$ws.Range("A43").Value = 49.83108108108108
$ws.Range("A44").Value = 49.1554054054054
$ws.Range("A45").Value = 48.141891891891895
$ws.Range("A46").Value = 48.310810810810814
$ws.Range("A47").Value = 48.141891891891895
$ws.Range("A48").Value = 51.35135135135135
$ws.Range("A49").Value = 49.32432432432432

# Row 50 is untouched by the commit (value stays 73.431734317343171).
